# ---------------------------------------------------------------------------
# Reproduces the commit "use git_upload_now 2019/08/12 週一 20:28:26.33":
#
#   1. The `_GoBack` bookmark (Word's "last edit location" marker) moves
#      from the end of the "...特發此函，以資證明。" paragraph to the very
#      start of the document (right before the first run of the first
#      paragraph, "推薦函").
#   2. The stale page-number field (" " + { PAGE \* MERGEFORMAT }) is
#      removed from the primary header, leaving only the
#      "臺中市立太平國民中學" text.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Reposition the _GoBack bookmark to the start of the document -------
# A bookmark collapsed exactly at position 0 gets auto-expanded to the
# first word by Bookmarks.Add, so we bracket a throwaway marker character,
# bookmark *that*, then delete the marker -- leaving a true zero-length
# bookmark sitting at the document start, exactly like Word leaves behind
# after an edit there.
$startOfDoc = $d.Range(0, 0)
$startOfDoc.InsertBefore("X")
$markerRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$d.Range(0, 1).Delete()

# --- 2. Strip the PAGE field (and its leading space run) from the header ---
$header = $d.Sections(1).Headers(1)

$pageField = $header.Range.Fields(1)
$pageField.Delete()

$leadingSpace = $header.Range
$leadingSpace.SetRange(0, 1)
$leadingSpace.Text = ""
